# Apply fixes based on feedback
# Target shape: slide 1, shape index 2 ("Group") - the reactive-strategy description textbox.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)

# ---------------------------------------------------------------------------
# 1) Reposition / resize the shape
# ---------------------------------------------------------------------------
# Target EMU values (from the authored diff):
#   off  x=688848  y=4052422
#   ext cx=4326212 cy=5639126
# The COM surface only exposes Left/Top/Width/Height in points (1 pt = 12700 EMU).
# The literals below were chosen so that the round-trip conversion lands on the
# exact target EMU values.
$sh.Left   = 54.23999977111817
$sh.Top    = 319.08833312988287
$sh.Width  = 340.6466217041016
$sh.Height = 444.0256805419922

# ---------------------------------------------------------------------------
# 2) Text edits
# ---------------------------------------------------------------------------
$tf2 = $sh.TextFrame2
$tr2 = $tf2.TextRange

function Replace-WholeText($rangeObj, $search, $newText) {
    # Replaces the first occurrence of $search (whole run or sub-run) with
    # $newText, re-reading the live text each time so indices stay valid.
    $full = $rangeObj.Text
    $idx = $full.IndexOf($search)
    if ($idx -lt 0) {
        Write-Host "NOT FOUND: [$search]"
        return -1
    }
    $sub = $rangeObj.Characters($idx + 1, $search.Length)
    $sub.Text = $newText
    return $idx
}

function Replace-FromOffset($rangeObj, $search, $newText, $startAt) {
    # Same as above, but starts searching from a given 0-based offset - used
    # to target a specific occurrence among several identical substrings.
    $full = $rangeObj.Text
    $idx = $full.IndexOf($search, $startAt)
    if ($idx -lt 0) {
        Write-Host "NOT FOUND: [$search] from $startAt"
        return -1
    }
    $sub = $rangeObj.Characters($idx + 1, $search.Length)
    $sub.Text = $newText
    return $idx
}

# -- "allows you to shortcut execution " -> split into three runs:
#    "allows you to shortcut " + "execution" + " "
# The combined text is unchanged; only the run boundaries differ (so a later
# independent formatting change could be applied to "execution" alone).
$full = $tr2.Text
$anchor = "allows you to shortcut execution "
$idx = $full.IndexOf($anchor)
$execStart = $idx + "allows you to shortcut ".Length
$execSub = $tr2.Characters($execStart + 1, "execution".Length)
$execSub.Text = "execution"

# -- "underlying resource detected as unhealthy" -> "... is detected ..."
Replace-WholeText $tr2 "underlying resource detected as unhealthy" "underlying resource is detected as unhealthy" | Out-Null

# -- "In " -> "In the " (three occurrences: Closed / Open / HalfOpen)
$pos = 0
$idx = Replace-FromOffset $tr2 "In " "In the " $pos
$pos = $idx + "In the ".Length
$idx = Replace-FromOffset $tr2 "In " "In the " $pos
$pos = $idx + "In the ".Length
$idx = Replace-FromOffset $tr2 "In " "In the " $pos

# -- Closed-state sentence: traffic -> invocations
Replace-WholeText $tr2 " state the circuit allows traffic to pass through and it monitors the failures." " state the circuit allows invocations to pass through and it monitors the failures." | Out-Null

# -- Open-state sentence: traffic -> invocations
Replace-WholeText $tr2 " state the circuit blocks traffic. " " state the circuit blocks invocations. " | Out-Null

# -- HalfOpen-state sentence: request -> invocation
Replace-WholeText $tr2 " state the circuit allows a single request to pass through as a probe." " state the circuit allows a single invocation to pass through as a probe." | Out-Null

# -- "if it was in Open state." -> "if it was in the Open state."
Replace-WholeText $tr2 " if it was in Open state." " if it was in the Open state." | Out-Null

# -- "The circuit shortcuts the execution with an " -> split into three runs:
#    "The circuit shortcuts " + "the execution " + "with an "
$full = $tr2.Text
$anchor2 = "The circuit shortcuts the execution with an "
$idx2 = $full.IndexOf($anchor2)
$midStart = $idx2 + "The circuit shortcuts ".Length
$midSub = $tr2.Characters($midStart + 1, "the execution ".Length)
$midSub.Text = "the execution "

# -- "if it was Isolated state." -> "if it was in the Isolated state."
Replace-WholeText $tr2 " if it was Isolated state." " if it was in the Isolated state." | Out-Null
